$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 21.24434529525963
$ws.Cells.Item(2, 3).Value = 15.98798310421177
$ws.Cells.Item(2, 4).Value = 3.901488561645881
$ws.Cells.Item(2, 5).Value = 29.33630424589722
$ws.Cells.Item(2, 6).Value = 21.19244645375618
$ws.Cells.Item(2, 7).Value = 2.068328951829761
$ws.Cells.Item(2, 8).Value = 3.23582982603127
$ws.Cells.Item(2, 9).Value = 3.346455671044311
$ws.Cells.Item(2, 16).Value = 11.7678710571247
$ws.Cells.Item(2, 17).Value = 16.6519096887458

$ws.Cells.Item(3, 2).Value = 19.83404491452785
$ws.Cells.Item(3, 3).Value = 15.06083541376267
$ws.Cells.Item(3, 4).Value = 3.796144748385419
$ws.Cells.Item(3, 5).Value = 27.45330212399358
$ws.Cells.Item(3, 6).Value = 20.3309374745291
$ws.Cells.Item(3, 7).Value = 2.07297968463036
$ws.Cells.Item(3, 8).Value = 2.982130108109396
$ws.Cells.Item(3, 9).Value = 3.168475253322779
$ws.Cells.Item(3, 16).Value = 11.85923295924826
$ws.Cells.Item(3, 17).Value = 16.1666207780694

$ws.Cells.Item(4, 2).Value = 18.91300103090603
$ws.Cells.Item(4, 3).Value = 14.46560476936423
$ws.Cells.Item(4, 4).Value = 3.729274292747749
$ws.Cells.Item(4, 5).Value = 26.23278851326196
$ws.Cells.Item(4, 6).Value = 19.792866582613
$ws.Cells.Item(4, 7).Value = 2.07592495572274
$ws.Cells.Item(4, 8).Value = 2.820997640627038
$ws.Cells.Item(4, 9).Value = 3.056221160496542
$ws.Cells.Item(4, 16).Value = 11.91791939148089
$ws.Cells.Item(4, 17).Value = 15.86810596097616

$ws.Cells.Item(5, 2).Value = 18.51737395584194
$ws.Cells.Item(5, 3).Value = 14.22662462568218
$ws.Cells.Item(5, 4).Value = 3.702647104088578
$ws.Cells.Item(5, 5).Value = 25.71893596989928
$ws.Cells.Item(5, 6).Value = 19.56228696292059
$ws.Cells.Item(5, 7).Value = 2.077155587032032
$ws.Cells.Item(5, 8).Value = 2.75378960930067
$ws.Cells.Item(5, 9).Value = 3.010004753504221
$ws.Cells.Item(5, 16).Value = 11.943404496475
$ws.Cells.Item(5, 17).Value = 15.73893713464746

$ws.Cells.Item(6, 2).Value = 18.44304980382731
$ws.Cells.Item(6, 3).Value = 14.198826641894
$ws.Cells.Item(6, 4).Value = 3.699614661890565
$ws.Cells.Item(6, 5).Value = 25.63228342230575
$ws.Cells.Item(6, 6).Value = 19.51249521288196
$ws.Cells.Item(6, 7).Value = 2.077370327591373
$ws.Cells.Item(6, 8).Value = 2.742335900622129
$ws.Cells.Item(6, 9).Value = 3.002823975041709
$ws.Cells.Item(6, 16).Value = 11.94882197734478
$ws.Cells.Item(6, 17).Value = 15.70830475881949

$ws.Cells.Item(7, 2).Value = 18.88696306560646
$ws.Cells.Item(7, 3).Value = 14.49524762205506
$ws.Cells.Item(7, 4).Value = 3.732781708952662
$ws.Cells.Item(7, 5).Value = 26.2249901650297
$ws.Cells.Item(7, 6).Value = 19.75907531799321
$ws.Cells.Item(7, 7).Value = 2.075965619517161
$ws.Cells.Item(7, 8).Value = 2.819531202845754
$ws.Cells.Item(7, 9).Value = 3.056754298870644
$ws.Cells.Item(7, 16).Value = 11.92148074146381
$ws.Cells.Item(7, 17).Value = 15.841509031171

$ws.Cells.Item(8, 2).Value = 20.74412174701266
$ws.Cells.Item(8, 3).Value = 15.71431410343474
$ws.Cells.Item(8, 4).Value = 3.870664877140791
$ws.Cells.Item(8, 5).Value = 28.69922864952596
$ws.Cells.Item(8, 6).Value = 20.85879796876124
$ws.Cells.Item(8, 7).Value = 2.069944485621055
$ws.Cells.Item(8, 8).Value = 3.14873931336443
$ws.Cells.Item(8, 9).Value = 3.286764084671139
$ws.Cells.Item(8, 16).Value = 11.80347848797396
$ws.Cells.Item(8, 17).Value = 16.45301393072154

$ws.Cells.Item(9, 2).Value = 23.9847716862231
$ws.Cells.Item(9, 3).Value = 17.84333741212275
$ws.Cells.Item(9, 4).Value = 4.116343728166592
$ws.Cells.Item(9, 5).Value = 33.04398820311501
$ws.Cells.Item(9, 6).Value = 22.98345331535855
$ws.Cells.Item(9, 7).Value = 2.058785775659262
$ws.Cells.Item(9, 8).Value = 3.753920340167167
$ws.Cells.Item(9, 9).Value = 3.714043418242903
$ws.Cells.Item(9, 16).Value = 11.58678339590977
$ws.Cells.Item(9, 17).Value = 17.68518224625051

$ws.Cells.Item(10, 2).Value = 26.03871136407511
$ws.Cells.Item(10, 3).Value = 19.23416383004324
$ws.Cells.Item(10, 4).Value = 4.311185115942086
$ws.Cells.Item(10, 5).Value = 35.05051666660535
$ws.Cells.Item(10, 6).Value = 24.32600997702151
$ws.Cells.Item(10, 7).Value = 2.051226890993335
$ws.Cells.Item(10, 8).Value = 4.137073897476883
$ws.Cells.Item(10, 9).Value = 4.007321852824707
$ws.Cells.Item(10, 16).Value = 11.46189746524951
$ws.Cells.Item(10, 17).Value = 18.44890189708534

$ws.Cells.Item(11, 2).Value = 26.39814353220931
$ws.Cells.Item(11, 3).Value = 19.39037959456375
$ws.Cells.Item(11, 4).Value = 4.617656838781644
$ws.Cells.Item(11, 5).Value = 28.4413378963673
$ws.Cells.Item(11, 6).Value = 23.65465203734853
$ws.Cells.Item(11, 7).Value = 2.0498429246093
$ws.Cells.Item(11, 8).Value = 4.492300152688726
$ws.Cells.Item(11, 9).Value = 4.071110765071338
$ws.Cells.Item(11, 16).Value = 11.58943319410997
$ws.Cells.Item(11, 17).Value = 17.73242486140492

$ws.Cells.Item(12, 2).Value = 26.31017654352046
$ws.Cells.Item(12, 3).Value = 19.20337017689211
$ws.Cells.Item(12, 4).Value = 4.838017044454928
$ws.Cells.Item(12, 5).Value = 22.48597269447575
$ws.Cells.Item(12, 6).Value = 22.83200781241673
$ws.Cells.Item(12, 7).Value = 2.050003634790813
$ws.Cells.Item(12, 8).Value = 5.297414652992155
$ws.Cells.Item(12, 9).Value = 4.068682360689743
$ws.Cells.Item(12, 16).Value = 11.72926721177756
$ws.Cells.Item(12, 17).Value = 16.98624232278696

$ws.Cells.Item(13, 2).Value = 25.84726470326447
$ws.Cells.Item(13, 3).Value = 18.77667581959529
$ws.Cells.Item(13, 4).Value = 5.01237998862362
$ws.Cells.Item(13, 5).Value = 16.56201812557444
$ws.Cells.Item(13, 6).Value = 21.78085996614549
$ws.Cells.Item(13, 7).Value = 2.051387570047776
$ws.Cells.Item(13, 8).Value = 6.317432706305129
$ws.Cells.Item(13, 9).Value = 4.017252459079821
$ws.Cells.Item(13, 16).Value = 11.88680967195127
$ws.Cells.Item(13, 17).Value = 16.12346427244232

$ws.Cells.Item(14, 2).Value = 25.34308545040505
$ws.Cells.Item(14, 3).Value = 18.37350547551109
$ws.Cells.Item(14, 4).Value = 5.116182715614461
$ws.Cells.Item(14, 5).Value = 12.49264091669018
$ws.Cells.Item(14, 6).Value = 20.92278044868427
$ws.Cells.Item(14, 7).Value = 2.052888961639587
$ws.Cells.Item(14, 8).Value = 7.116731744846066
$ws.Cells.Item(14, 9).Value = 3.960504999053573
$ws.Cells.Item(14, 16).Value = 12.00576682982438
$ws.Cells.Item(14, 17).Value = 15.45461220619254

$ws.Cells.Item(15, 2).Value = 25.13125716500807
$ws.Cells.Item(15, 3).Value = 18.23190355643588
$ws.Cells.Item(15, 4).Value = 5.132674718159113
$ws.Cells.Item(15, 5).Value = 11.51593739001252
$ws.Cells.Item(15, 6).Value = 20.65438187367861
$ws.Cells.Item(15, 7).Value = 2.053551800275636
$ws.Cells.Item(15, 8).Value = 7.300747388952847
$ws.Cells.Item(15, 9).Value = 3.9361343313207
$ws.Cells.Item(15, 16).Value = 12.0386306197751
$ws.Cells.Item(15, 17).Value = 15.25954650730678

$ws.Cells.Item(16, 2).Value = 24.34195033592232
$ws.Cells.Item(16, 3).Value = 17.73640174565618
$ws.Cells.Item(16, 4).Value = 5.022720057685066
$ws.Cells.Item(16, 5).Value = 11.39309845988134
$ws.Cells.Item(16, 6).Value = 20.22206244717581
$ws.Cells.Item(16, 7).Value = 2.056465447626719
$ws.Cells.Item(16, 8).Value = 7.02225748233182
$ws.Cells.Item(16, 9).Value = 3.824142457375225
$ws.Cells.Item(16, 16).Value = 12.04649673184934
$ws.Cells.Item(16, 17).Value = 15.06049289056866

$ws.Cells.Item(17, 2).Value = 24.00371943281733
$ws.Cells.Item(17, 3).Value = 17.57351854959033
$ws.Cells.Item(17, 4).Value = 4.880251754739056
$ws.Cells.Item(17, 5).Value = 13.4979010985307
$ws.Cells.Item(17, 6).Value = 20.36296271496995
$ws.Cells.Item(17, 7).Value = 2.057930077344666
$ws.Cells.Item(17, 8).Value = 6.332660876888924
$ws.Cells.Item(17, 9).Value = 3.767659763510585
$ws.Cells.Item(17, 16).Value = 11.99182856646824
$ws.Cells.Item(17, 17).Value = 15.27309957087339

$ws.Cells.Item(18, 2).Value = 24.044348790412
$ws.Cells.Item(18, 3).Value = 17.66693945972051
$ws.Cells.Item(18, 4).Value = 4.693155499359199
$ws.Cells.Item(18, 5).Value = 18.0617373087299
$ws.Cells.Item(18, 6).Value = 21.03231258682181
$ws.Cells.Item(18, 7).Value = 2.058178860574236
$ws.Cells.Item(18, 8).Value = 5.299013535784304
$ws.Cells.Item(18, 9).Value = 3.754025291857811
$ws.Cells.Item(18, 16).Value = 11.87671840079358
$ws.Cells.Item(18, 17).Value = 15.87749421297165

$ws.Cells.Item(19, 2).Value = 24.35879215894935
$ws.Cells.Item(19, 3).Value = 18.01543152623328
$ws.Cells.Item(19, 4).Value = 4.497482009348721
$ws.Cells.Item(19, 5).Value = 24.34212166184587
$ws.Cells.Item(19, 6).Value = 22.01514394023083
$ws.Cells.Item(19, 7).Value = 2.057304220750735
$ws.Cells.Item(19, 8).Value = 4.323344623511946
$ws.Cells.Item(19, 9).Value = 3.785884445289801
$ws.Cells.Item(19, 16).Value = 11.73771417547752
$ws.Cells.Item(19, 17).Value = 16.70409758324919

$ws.Cells.Item(20, 2).Value = 25.46361636489998
$ws.Cells.Item(20, 3).Value = 18.95625835189501
$ws.Cells.Item(20, 4).Value = 4.273331612968068
$ws.Cells.Item(20, 5).Value = 34.4935351717964
$ws.Cells.Item(20, 6).Value = 23.8909754692182
$ws.Cells.Item(20, 7).Value = 2.053260626197614
$ws.Cells.Item(20, 8).Value = 4.03386592142914
$ws.Cells.Item(20, 9).Value = 3.935342292896853
$ws.Cells.Item(20, 16).Value = 11.50748466950926
$ws.Cells.Item(20, 17).Value = 18.17589669579615

$ws.Cells.Item(21, 2).Value = 27.06302406206178
$ws.Cells.Item(21, 3).Value = 20.05080150246729
$ws.Cells.Item(21, 4).Value = 4.378876214376258
$ws.Cells.Item(21, 5).Value = 37.30861685159484
$ws.Cells.Item(21, 6).Value = 25.14104151860499
$ws.Cells.Item(21, 7).Value = 2.047180558833289
$ws.Cells.Item(21, 8).Value = 4.380770126744784
$ws.Cells.Item(21, 9).Value = 4.169007924582534
$ws.Cells.Item(21, 16).Value = 11.38581223215721
$ws.Cells.Item(21, 17).Value = 18.95822593577839

$ws.Cells.Item(22, 2).Value = 28.05083845722103
$ws.Cells.Item(22, 3).Value = 20.6820002794959
$ws.Cells.Item(22, 4).Value = 4.453456234324586
$ws.Cells.Item(22, 5).Value = 38.6401310496351
$ws.Cells.Item(22, 6).Value = 25.90144730192458
$ws.Cells.Item(22, 7).Value = 2.043365079700162
$ws.Cells.Item(22, 8).Value = 4.584481656539372
$ws.Cells.Item(22, 9).Value = 4.313709910463947
$ws.Cells.Item(22, 16).Value = 11.31264952474652
$ws.Cells.Item(22, 17).Value = 19.43138575637101

$ws.Cells.Item(23, 2).Value = 27.54557798507469
$ws.Cells.Item(23, 3).Value = 20.3181632356919
$ws.Cells.Item(23, 4).Value = 4.409673853925127
$ws.Cells.Item(23, 5).Value = 37.93528199551817
$ws.Cells.Item(23, 6).Value = 25.52562432809528
$ws.Cells.Item(23, 7).Value = 2.045373237940637
$ws.Cells.Item(23, 8).Value = 4.476635571788233
$ws.Cells.Item(23, 9).Value = 4.234884081058053
$ws.Cells.Item(23, 16).Value = 11.34622501537988
$ws.Cells.Item(23, 17).Value = 19.20377684041285

$ws.Cells.Item(24, 2).Value = 25.51262758496888
$ws.Cells.Item(24, 3).Value = 18.93701595186494
$ws.Cells.Item(24, 4).Value = 4.245717152561948
$ws.Cells.Item(24, 5).Value = 35.15552401392027
$ws.Cells.Item(24, 6).Value = 24.03099774490324
$ws.Cells.Item(24, 7).Value = 2.053131939213595
$ws.Cells.Item(24, 8).Value = 4.059945008673057
$ws.Cells.Item(24, 9).Value = 3.934995918583402
$ws.Cells.Item(24, 16).Value = 11.48737804446741
$ws.Cells.Item(24, 17).Value = 18.29734762139054

$ws.Cells.Item(25, 2).Value = 23.12287659847034
$ws.Cells.Item(25, 3).Value = 17.3406112688562
$ws.Cells.Item(25, 4).Value = 4.058593364205172
$ws.Cells.Item(25, 5).Value = 31.92129429900232
$ws.Cells.Item(25, 6).Value = 22.36976730959418
$ws.Cells.Item(25, 7).Value = 2.061772739636632
$ws.Cells.Item(25, 8).Value = 3.593540973059603
$ws.Cells.Item(25, 9).Value = 3.602901949824975
$ws.Cells.Item(25, 16).Value = 11.65009461266547
$ws.Cells.Item(25, 17).Value = 17.31236120644557
